$d = $word.ActiveDocument

# 1) Word choice fix: "суперполиномиальной" -> "сверх-полиномиальной"
$d.Content.Find.Execute(
    "суперполиномиальной", $false, $false, $false, $false, $false,
    $true, 1, $false, "сверх-полиномиальной", 2) | Out-Null

# 2) Typo fix: "одностороенней" -> "односторонней"
$d.Content.Find.Execute(
    "одностороенней", $false, $false, $false, $false, $false,
    $true, 1, $false, "односторонней", 2) | Out-Null

# 3) Remove duplicated space: "стойкая  к коллизиям" -> "стойкая к коллизиям"
$d.Content.Find.Execute(
    "стойкая  к коллизиям хэш-функция. Доказать от противного",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "стойкая к коллизиям хэш-функция. Доказать от противного", 2) | Out-Null

# 4) Move the "_GoBack" bookmark from the end of the document to right
#    after the just-fixed "односторонней" word (before " хэш-функции").
$anchor = $d.Content
$anchor.Find.Execute(
    "в модели односторонней хэш-функции", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$markPos = $anchor.End - " хэш-функции".Length
$bookmarkRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
